$d = $word.ActiveDocument

# Helper: given a character offset into the document, return the
# Paragraph object whose Range actually contains it. (Range.Paragraphs
# on an arbitrary sub-range does not reliably expand to the enclosing
# paragraph in this host, so we walk the document's Paragraphs
# collection instead, using a half-open [Start, End) test so boundary
# offsets resolve to the paragraph that *starts* there.)
function Get-ParagraphContainingOffset($doc, $offset) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Start -le $offset -and $offset -lt $p.Range.End) {
            return $p
        }
    }
    return $doc.Paragraphs.Item($doc.Paragraphs.Count)
}

# Find the paragraph that says "Ver no Jupiter Salvar em pdf Salvar em docx".
$startRng = $d.Content
$startRng.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $startRng.Find.Found) {
    throw "Could not find the 'Ver no Jupiter...' paragraph"
}
$startPara = Get-ParagraphContainingOffset $d $startRng.Start

# Find the paragraph with the copyright / footer notice.
$endRng = $d.Content
$endRng.Find.Execute("Contact: luizeleno@usp.br", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $endRng.Find.Found) {
    throw "Could not find the copyright paragraph"
}
$copyPara = Get-ParagraphContainingOffset $d $endRng.Start

# The (now redundant) blank paragraph that used to separate the footer
# block from the following page-break paragraph is removed too.
$blankPara = $copyPara.Next()

$deleteRange = $d.Range($startPara.Range.Start, $blankPara.Range.End)
$deleteRange.Delete()
